# PASADOR CERROJO.xlsx — update date and prices on Hoja1
#
# Changes (per commit diff):
#   A1  (date, numFmt 14)      45406  -> 45436   (2024-04-24 -> 2024-05-24)
#   D29 (price, "PI-1500")     185.28 -> 364.992
#   D30 (price, "PI-1501")     261.067 -> 514.29
#
# Only the cell values change; everything else (layout, merges, styles)
# stays as-is, so we touch nothing but these three cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436
$ws.Range("D29").Value = 364.992
$ws.Range("D30").Value = 514.29
